# Update the Production_Predictions worksheet data.
# The underlying data window shifted back by two days (and the trailing
# duplicated/garbled rows were dropped), so every data row's Date/Interval/
# Prediction values are rewritten in place, and the now-unused last row
# (row 74) is removed so the sheet ends at row 73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, DateSerial, Interval, Prediction
$data = @(
    @(2, 45342, 0, 0.004000000189989805),
    @(3, 45342, 1, 0.004000000189989805),
    @(4, 45342, 2, 0.004000000189989805),
    @(5, 45342, 3, 0.004000000189989805),
    @(6, 45342, 4, 0.004000000189989805),
    @(7, 45342, 5, 0.004000000189989805),
    @(8, 45342, 6, 0.004000000189989805),
    @(9, 45342, 7, 0.00800000037997961),
    @(10, 45342, 8, 0.0260000005364418),
    @(11, 45342, 9, 0.05400000140070915),
    @(12, 45342, 10, 0.07299999892711639),
    @(13, 45342, 11, 0.08600000292062759),
    @(14, 45342, 12, 0.0820000022649765),
    @(15, 45342, 13, 0.08100000023841858),
    @(16, 45342, 14, 0.06800000369548798),
    @(17, 45342, 15, 0.05499999970197678),
    @(18, 45342, 16, 0.0260000005364418),
    @(19, 45342, 17, 0.008999999612569809),
    @(20, 45342, 18, 0.004000000189989805),
    @(21, 45342, 19, 0.004000000189989805),
    @(22, 45342, 20, 0.004000000189989805),
    @(23, 45342, 21, 0.004000000189989805),
    @(24, 45342, 22, 0.004000000189989805),
    @(25, 45342, 23, 0.004000000189989805),
    @(26, 45343, 0, 0.004000000189989805),
    @(27, 45343, 1, 0.004000000189989805),
    @(28, 45343, 2, 0.004000000189989805),
    @(29, 45343, 3, 0.004000000189989805),
    @(30, 45343, 4, 0.004000000189989805),
    @(31, 45343, 5, 0.004000000189989805),
    @(32, 45343, 6, 0.004000000189989805),
    @(33, 45343, 7, 0.00800000037997961),
    @(34, 45343, 8, 0.0390000008046627),
    @(35, 45343, 9, 0.06300000101327896),
    @(36, 45343, 10, 0.08799999952316284),
    @(37, 45343, 11, 0.1040000021457672),
    @(38, 45343, 12, 0.09099999815225601),
    @(39, 45343, 13, 0.0949999988079071),
    @(40, 45343, 14, 0.08100000023841858),
    @(41, 45343, 15, 0.07199999690055847),
    @(42, 45343, 16, 0.03599999845027924),
    @(43, 45343, 17, 0.01099999994039536),
    @(44, 45343, 18, 0.004000000189989805),
    @(45, 45343, 19, 0.004000000189989805),
    @(46, 45343, 20, 0.004000000189989805),
    @(47, 45343, 21, 0.004000000189989805),
    @(48, 45343, 22, 0.004000000189989805),
    @(49, 45343, 23, 0.004000000189989805),
    @(50, 45344, 0, 0.004000000189989805),
    @(51, 45344, 1, 0.004000000189989805),
    @(52, 45344, 2, 0.004000000189989805),
    @(53, 45344, 3, 0.004000000189989805),
    @(54, 45344, 4, 0.004000000189989805),
    @(55, 45344, 5, 0.004000000189989805),
    @(56, 45344, 6, 0.004000000189989805),
    @(57, 45344, 7, 0.00800000037997961),
    @(58, 45344, 8, 0.03799999877810478),
    @(59, 45344, 9, 0.06300000101327896),
    @(60, 45344, 10, 0.08699999749660492),
    @(61, 45344, 11, 0.1230000033974648),
    @(62, 45344, 12, 0.125),
    @(63, 45344, 13, 0.1270000040531158),
    @(64, 45344, 14, 0.1169999986886978),
    @(65, 45344, 15, 0.1220000013709068),
    @(66, 45344, 16, 0.07100000232458115),
    @(67, 45344, 17, 0.01099999994039536),
    @(68, 45344, 18, 0.004000000189989805),
    @(69, 45344, 19, 0.004000000189989805),
    @(70, 45344, 20, 0.004000000189989805),
    @(71, 45344, 21, 0.004000000189989805),
    @(72, 45344, 22, 0.004000000189989805),
    @(73, 45344, 23, 0.004000000189989805)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
}

# The data table now only spans through row 73; remove the old trailing row 74.
$ws.Rows(74).Delete()
